$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2047244094488189
$ws.Range("C2").Value = 0.531496062992126
$ws.Range("J2").Value = 0.01181102362204724
$ws.Range("P2").Value = 0.1535433070866142
$ws.Range("S2").Value = 0.09842519685039371
$ws.Range("B3").Value = 0.02054794520547945
$ws.Range("C3").Value = 0.0547945205479452
$ws.Range("J3").Value = 0.02054794520547945
$ws.Range("P3").Value = 0.8013698630136986
$ws.Range("S3").Value = 0.1027397260273973
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.303030303030303
$ws.Range("B6").Value = 0.03867403314917127
$ws.Range("D6").Value = 0.01657458563535912
$ws.Range("F6").Value = 0.06629834254143646
$ws.Range("J6").Value = 0.2375690607734807
$ws.Range("O6").Value = 0.01657458563535912
$ws.Range("Q6").Value = 0.1657458563535912
$ws.Range("R6").Value = 0.0718232044198895
$ws.Range("S6").Value = 0.3867403314917127
$ws.Range("B7").Value = 0.09239130434782608
$ws.Range("D7").Value = 0.0108695652173913
$ws.Range("F7").Value = 0.07065217391304347
$ws.Range("J7").Value = 0.04891304347826087
$ws.Range("O7").Value = 0.03260869565217391
$ws.Range("Q7").Value = 0.1521739130434783
$ws.Range("R7").Value = 0.125
$ws.Range("S7").Value = 0.4673913043478261
$ws.Range("B8").Value = 0.09534368070953436
$ws.Range("D8").Value = 0.01995565410199556
$ws.Range("F8").Value = 0.05543237250554324
$ws.Range("J8").Value = 0.1042128603104213
$ws.Range("O8").Value = 0.03104212860310421
$ws.Range("Q8").Value = 0.1862527716186253
$ws.Range("R8").Value = 0.1152993348115299
$ws.Range("S8").Value = 0.3924611973392461
$ws.Range("B9").Value = 0.06696428571428571
$ws.Range("D9").Value = 0.02232142857142857
$ws.Range("F9").Value = 0.04017857142857143
$ws.Range("J9").Value = 0.08035714285714286
$ws.Range("O9").Value = 0.04464285714285714
$ws.Range("Q9").Value = 0.1741071428571428
$ws.Range("R9").Value = 0.1294642857142857
$ws.Range("S9").Value = 0.4419642857142857
$ws.Range("B10").Value = 0.1028880866425993
$ws.Range("D10").Value = 0.01353790613718412
$ws.Range("F10").Value = 0.05415162454873646
$ws.Range("J10").Value = 0.07851985559566788
$ws.Range("O10").Value = 0.01624548736462094
$ws.Range("Q10").Value = 0.2129963898916968
$ws.Range("R10").Value = 0.1010830324909747
$ws.Range("S10").Value = 0.4205776173285198
$ws.Range("G11").Value = 0.1314878892733564
$ws.Range("J11").Value = 0.08304498269896193
$ws.Range("K11").Value = 0.1764705882352941
$ws.Range("L11").Value = 0.5709342560553633
$ws.Range("S11").Value = 0.03806228373702422
$ws.Range("G12").Value = 0.7455621301775148
$ws.Range("J12").Value = 0.1715976331360947
$ws.Range("K12").Value = 0.005917159763313609
$ws.Range("L12").Value = 0.02958579881656805
$ws.Range("S12").Value = 0.04733727810650887
$ws.Range("F13").Value = 0.02439024390243903
$ws.Range("G13").Value = 0.6585365853658537
$ws.Range("J13").Value = 0.2195121951219512
$ws.Range("S13").Value = 0.0975609756097561
$ws.Range("F15").Value = 0.0179372197309417
$ws.Range("H15").Value = 0.1659192825112108
$ws.Range("I15").Value = 0.07623318385650224
$ws.Range("J15").Value = 0.2869955156950673
$ws.Range("K15").Value = 0.07174887892376682
$ws.Range("M15").Value = 0.008968609865470852
$ws.Range("O15").Value = 0.04484304932735426
$ws.Range("S15").Value = 0.3273542600896861
$ws.Range("F16").Value = 0.02352941176470588
$ws.Range("H16").Value = 0.1823529411764706
$ws.Range("I16").Value = 0.07058823529411765
$ws.Range("J16").Value = 0.4235294117647059
$ws.Range("K16").Value = 0.08823529411764706
$ws.Range("M16").Value = 0.02352941176470588
$ws.Range("N16").Value = 0.005882352941176471
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.1235294117647059
$ws.Range("F17").Value = 0.02147971360381861
$ws.Range("H17").Value = 0.1813842482100239
$ws.Range("I17").Value = 0.1002386634844869
$ws.Range("J17").Value = 0.431980906921241
$ws.Range("K17").Value = 0.07875894988066826
$ws.Range("M17").Value = 0.01431980906921241
$ws.Range("N17").Value = 0.002386634844868735
$ws.Range("O17").Value = 0.06682577565632458
$ws.Range("S17").Value = 0.1026252983293556
$ws.Range("F18").Value = 0.01293103448275862
$ws.Range("H18").Value = 0.146551724137931
$ws.Range("I18").Value = 0.1163793103448276
$ws.Range("J18").Value = 0.3879310344827586
$ws.Range("K18").Value = 0.1120689655172414
$ws.Range("M18").Value = 0.02586206896551724
$ws.Range("O18").Value = 0.04310344827586207
$ws.Range("S18").Value = 0.1551724137931035
$ws.Range("F19").Value = 0.01059803179409538
$ws.Range("H19").Value = 0.2096896290688872
$ws.Range("I19").Value = 0.09538228614685844
$ws.Range("J19").Value = 0.3414080242240727
$ws.Range("K19").Value = 0.1044663133989402
$ws.Range("M19").Value = 0.01968205904617714
$ws.Range("O19").Value = 0.06661619984859955
$ws.Range("S19").Value = 0.1521574564723694
